$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Year" header label to A1 (B1 already holds "Population")
$ws.Range("A1").Value = "Year"

# Move the active selection to C4, matching the saved selection state
$ws.Range("C4").Select()
